# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45192 to serial date 45202.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 116; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45202
}
